# Apply scheduled market-data updates to the Exodus Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 98995.664
$ws.Range("J133").Value = 98995.664
$ws.Range("L133").Value = 98995.664
$ws.Range("N133").Value = -109115.664
$ws.Range("H134").Value = 93082.5
$ws.Range("J134").Value = 93082.5
$ws.Range("L134").Value = 93082.5
$ws.Range("N134").Value = -103222.5
$ws.Range("H136").Value = 77977.14
$ws.Range("J136").Value = 77977.14
$ws.Range("L136").Value = 77977.14
$ws.Range("N136").Value = -88177.14
$ws.Range("H137").Value = 394779.9
$ws.Range("J137").Value = 633768.6
$ws.Range("L137").Value = 1901305.8
$ws.Range("N137").Value = -1906405.8
$ws.Range("H139").Value = 70168.44500000001
$ws.Range("J139").Value = 70168.44500000001
$ws.Range("L139").Value = 70168.44500000001
$ws.Range("N139").Value = -80448.44500000001
$ws.Range("H140").Value = 85716.37
$ws.Range("J140").Value = 91217.10000000001
$ws.Range("L140").Value = 91217.10000000001
$ws.Range("N140").Value = -101577.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 5750
$ws.Range("I38").Value = 2750
$ws.Range("J38").Value = 8750
$ws.Range("K38").Value = 2750
$ws.Range("L38").Value = 8750
$ws.Range("M38").Value = -2283
$ws.Range("N38").Value = -9684
$ws.Range("H39").Value = 10400.2
$ws.Range("I39").Value = 10400.2
$ws.Range("K39").Value = 10400.2
$ws.Range("M39").Value = -9880.200000000001
$ws.Range("H61").Value = 73559.64
$ws.Range("I61").Value = 1683.7
$ws.Range("K61").Value = 1683.7
$ws.Range("M61").Value = -1471.7
$ws.Range("H88").Value = 1635.3846
$ws.Range("J88").Value = 2093.25
$ws.Range("L88").Value = 2093.25
$ws.Range("N88").Value = -2905.25
$ws.Range("H91").Value = 1635.3846
$ws.Range("J91").Value = 2093.25
$ws.Range("L91").Value = 2093.25
$ws.Range("N91").Value = -4901.25
$ws.Range("H102").Value = 94081.836
$ws.Range("I102").Value = 102572.1
$ws.Range("J102").Value = 51630.5
$ws.Range("K102").Value = 102572.1
$ws.Range("L102").Value = 51630.5
$ws.Range("M102").Value = -100950.1
$ws.Range("N102").Value = -54874.5
$ws.Range("H136").Value = 73559.64
$ws.Range("I136").Value = 1683.7
$ws.Range("K136").Value = 5051.1
$ws.Range("M136").Value = -2501.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 168957.5
$ws.Range("I105").Value = 500450
$ws.Range("J105").Value = 3211.25
$ws.Range("K105").Value = 500450
$ws.Range("L105").Value = 3211.25
$ws.Range("M105").Value = -498703
$ws.Range("N105").Value = -6705.25
$ws.Range("H107").Value = 2721.158
$ws.Range("I107").Value = 2621.3845
$ws.Range("J107").Value = 2937.3333
$ws.Range("K107").Value = 2621.3845
$ws.Range("L107").Value = 2937.3333
$ws.Range("M107").Value = -701.3845000000001
$ws.Range("N107").Value = -6777.3333
$ws.Range("H132").Value = 51194.445
$ws.Range("J132").Value = 51194.445
$ws.Range("L132").Value = 51194.445
$ws.Range("N132").Value = -61314.445
$ws.Range("H135").Value = 100440
$ws.Range("J135").Value = 100440
$ws.Range("L135").Value = 100440
$ws.Range("N135").Value = -110580
$ws.Range("H138").Value = 80711.42999999999
$ws.Range("J138").Value = 80711.42999999999
$ws.Range("L138").Value = 80711.42999999999
$ws.Range("N138").Value = -90991.42999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3371.5476
$ws.Range("I58").Value = 3444.238
$ws.Range("K58").Value = 3444.238
$ws.Range("M58").Value = -3241.238
$ws.Range("H136").Value = 3371.5476
$ws.Range("I136").Value = 3444.238
$ws.Range("K136").Value = 10332.714
$ws.Range("M136").Value = -7782.714
$ws.Range("H138").Value = 89996
$ws.Range("J138").Value = 89996
$ws.Range("L138").Value = 89996
$ws.Range("N138").Value = -100276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 905.6
$ws.Range("I5").Value = 1081.8
$ws.Range("J5").Value = 817.5
$ws.Range("K5").Value = 3245.4
$ws.Range("L5").Value = 2452.5
$ws.Range("M5").Value = -3133.4
$ws.Range("N5").Value = -2676.5
$ws.Range("H60").Value = 215.20589
$ws.Range("I60").Value = 58.083332
$ws.Range("K60").Value = 174.249996
$ws.Range("M60").Value = 76.75000399999999
$ws.Range("H124").Value = 4777
$ws.Range("I124").Value = 4777
$ws.Range("K124").Value = 14331
$ws.Range("M124").Value = -9421
$ws.Range("H125").Value = 7110.75
$ws.Range("J125").Value = 15000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H129").Value = 2416.0908
$ws.Range("I129").Value = 753.8333
$ws.Range("J129").Value = 4410.8
$ws.Range("K129").Value = 2261.4999
$ws.Range("L129").Value = 13232.4
$ws.Range("M129").Value = 2738.5001
$ws.Range("N129").Value = -23232.4
$ws.Range("H131").Value = 1718.5
$ws.Range("J131").Value = 2149.6
$ws.Range("L131").Value = 6448.799999999999
$ws.Range("N131").Value = -16528.8
$ws.Range("H132").Value = 8412.643
$ws.Range("J132").Value = 8412.643
$ws.Range("L132").Value = 75713.787
$ws.Range("N132").Value = -80773.787
$ws.Range("H133").Value = 6789.2856
$ws.Range("I133").Value = 5505
$ws.Range("K133").Value = 16515
$ws.Range("M133").Value = -11455
$ws.Range("H135").Value = 905.6
$ws.Range("I135").Value = 1081.8
$ws.Range("J135").Value = 817.5
$ws.Range("K135").Value = 9736.199999999999
$ws.Range("L135").Value = 7357.5
$ws.Range("M135").Value = -7201.199999999999
$ws.Range("N135").Value = -12427.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 5679.6665
$ws.Range("J47").Value = 5875.6
$ws.Range("L47").Value = 5875.6
$ws.Range("N47").Value = -7011.6
$ws.Range("H80").Value = 41016.31
$ws.Range("I80").Value = 2340.4285
$ws.Range("J80").Value = 86138.164
$ws.Range("K80").Value = 2340.4285
$ws.Range("L80").Value = 86138.164
$ws.Range("M80").Value = -1342.4285
$ws.Range("N80").Value = -88134.164
$ws.Range("H83").Value = 41016.31
$ws.Range("I83").Value = 2340.4285
$ws.Range("J83").Value = 86138.164
$ws.Range("K83").Value = 11702.1425
$ws.Range("L83").Value = 430690.82
$ws.Range("M83").Value = -6710.1425
$ws.Range("N83").Value = -440674.82
$ws.Range("H93").Value = 18110
$ws.Range("J93").Value = 18110
$ws.Range("L93").Value = 18110
$ws.Range("N93").Value = -21854
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820
$ws.Range("H126").Value = 2967.1428
$ws.Range("I126").Value = 2263.647
$ws.Range("J126").Value = 4054.3635
$ws.Range("K126").Value = 6790.941
$ws.Range("L126").Value = 12163.0905
$ws.Range("M126").Value = -4320.941
$ws.Range("N126").Value = -17103.0905
$ws.Range("H132").Value = 7353.5
$ws.Range("I132").Value = 2544.9
$ws.Range("K132").Value = 7634.700000000001
$ws.Range("M132").Value = -5104.700000000001
$ws.Range("H135").Value = 98664.39999999999
$ws.Range("J135").Value = 98664.39999999999
$ws.Range("L135").Value = 98664.39999999999
$ws.Range("N135").Value = -108804.4
$ws.Range("H140").Value = 81067
$ws.Range("J140").Value = 92626.75
$ws.Range("L140").Value = 92626.75
$ws.Range("N140").Value = -102986.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5305.9375
$ws.Range("I132").Value = 3772.3635
$ws.Range("J132").Value = 8679.799999999999
$ws.Range("K132").Value = 11317.0905
$ws.Range("L132").Value = 26039.4
$ws.Range("M132").Value = -8787.0905
$ws.Range("N132").Value = -31099.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6532.4287
$ws.Range("I81").Value = 6496.7144
$ws.Range("K81").Value = 12993.4288
$ws.Range("M81").Value = -11932.4288
$ws.Range("H84").Value = 6532.4287
$ws.Range("I84").Value = 6496.7144
$ws.Range("K84").Value = 64967.144
$ws.Range("M84").Value = -59663.144
$ws.Range("H132").Value = 2077.8718
$ws.Range("I132").Value = 1551.1333
$ws.Range("K132").Value = 4653.3999
$ws.Range("M132").Value = -2123.3999
$ws.Range("H133").Value = 149999
$ws.Range("J133").Value = 149999
$ws.Range("L133").Value = 149999
$ws.Range("N133").Value = -160119
$ws.Range("H136").Value = 3383.7896
$ws.Range("I136").Value = 3268.25
$ws.Range("K136").Value = 9804.75
$ws.Range("M136").Value = -7254.75
$ws.Range("H140").Value = 126663.664
$ws.Range("J140").Value = 126663.664
$ws.Range("L140").Value = 126663.664
$ws.Range("N140").Value = -137023.664
